# Update the "Pais" (countries) sheet with the latest COVID-19 snapshot:
#  - refresh the "Datos actualizados..." timestamp
#  - update case/death/recovered figures for countries whose numbers changed
#  - a handful of countries changed ranking position, so their name (column A)
#    is rewritten in place to reflect the new row they occupy while the
#    updated statistics for that rank are written alongside it
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Marzo de 2020 a las 11:46"
# Row 6
$ws.Cells.Item(6, 1).Value = "España"
$ws.Cells.Item(6, 2).Value = 19980
$ws.Cells.Item(6, 3).Value = 1903
$ws.Cells.Item(6, 4).Value = 1588
$ws.Cells.Item(6, 5).Value = 17390
$ws.Cells.Item(6, 6).Value = 939
$ws.Cells.Item(6, 7).Value = 171
$ws.Cells.Item(6, 8).Value = 1002
# Row 7
$ws.Cells.Item(7, 1).Value = "Iran"
$ws.Cells.Item(7, 2).Value = 18407
$ws.Cells.Item(7, 4).Value = 5979
$ws.Cells.Item(7, 5).Value = 11144
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 8).Value = 1284
# Row 8
$ws.Cells.Item(8, 2).Value = 16626
$ws.Cells.Item(8, 3).Value = 1306
$ws.Cells.Item(8, 5).Value = 16467
# Row 15
$ws.Cells.Item(15, 4).Value = 204
$ws.Cells.Item(15, 5).Value = 2016
$ws.Cells.Item(15, 6).Value = 164
# Row 35
$ws.Cells.Item(35, 1).Value = "Polonia"
$ws.Cells.Item(35, 2).Value = 378
$ws.Cells.Item(35, 3).Value = 23
$ws.Cells.Item(35, 4).Value = 13
$ws.Cells.Item(35, 5).Value = 359
$ws.Cells.Item(35, 6).Value = 3
$ws.Cells.Item(35, 7).Value = 1
$ws.Cells.Item(35, 8).Value = 6
# Row 36
$ws.Cells.Item(36, 1).Value = "Indonesia"
$ws.Cells.Item(36, 2).Value = 369
$ws.Cells.Item(36, 3).Value = 60
$ws.Cells.Item(36, 4).Value = 17
$ws.Cells.Item(36, 5).Value = 320
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 7
$ws.Cells.Item(36, 8).Value = 32
# Row 43
$ws.Cells.Item(43, 1).Value = "Barein"
$ws.Cells.Item(43, 2).Value = 284
$ws.Cells.Item(43, 3).Value = 5
$ws.Cells.Item(43, 4).Value = 110
$ws.Cells.Item(43, 5).Value = 173
$ws.Cells.Item(43, 6).Value = 4
$ws.Cells.Item(43, 8).Value = 1
# Row 44
$ws.Cells.Item(44, 1).Value = "Estonia"
$ws.Cells.Item(44, 2).Value = 283
$ws.Cells.Item(44, 3).Value = 16
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(44, 5).Value = 282
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 8).Value = 0
# Row 75
$ws.Cells.Item(75, 5).Value = 75
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = 3
# Row 77
$ws.Cells.Item(77, 2).Value = 75
$ws.Cells.Item(77, 3).Value = 1
$ws.Cells.Item(77, 5).Value = 74
$ws.Cells.Item(77, 6).Value = 2
# Row 84
$ws.Cells.Item(84, 1).Value = "Sri Lanka"
$ws.Cells.Item(84, 3).Value = 6
$ws.Cells.Item(84, 4).Value = 3
$ws.Cells.Item(84, 5).Value = 63
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
# Row 85
$ws.Cells.Item(85, 1).Value = "Marruecos"
$ws.Cells.Item(85, 2).Value = 66
$ws.Cells.Item(85, 3).Value = 3
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 5).Value = 61
$ws.Cells.Item(85, 6).Value = 1
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 3
# Row 87
$ws.Cells.Item(87, 1).Value = "Tunez"
$ws.Cells.Item(87, 2).Value = 54
$ws.Cells.Item(87, 3).Value = 15
$ws.Cells.Item(87, 4).Value = 1
$ws.Cells.Item(87, 5).Value = 52
$ws.Cells.Item(87, 6).Value = 2
$ws.Cells.Item(87, 8).Value = 1
# Row 88
$ws.Cells.Item(88, 1).Value = "Malta"
$ws.Cells.Item(88, 2).Value = 53
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 2
$ws.Cells.Item(88, 5).Value = 51
# Row 89
$ws.Cells.Item(89, 1).Value = "Kazajistan"
$ws.Cells.Item(89, 3).Value = 5
$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(89, 5).Value = 49
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 8).Value = 0
# Row 90
$ws.Cells.Item(90, 1).Value = "Moldavia"
$ws.Cells.Item(90, 2).Value = 49
$ws.Cells.Item(90, 6).Value = 3
$ws.Cells.Item(90, 8).Value = 1
# Row 91
$ws.Cells.Item(91, 1).Value = "Lituania"
$ws.Cells.Item(91, 4).Value = 1
$ws.Cells.Item(91, 5).Value = 47
$ws.Cells.Item(91, 6).Value = 1
# Row 92
$ws.Cells.Item(92, 1).Value = "Oman"
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 13
$ws.Cells.Item(92, 5).Value = 35
# Row 93
$ws.Cells.Item(93, 1).Value = "Estado de Palestina"
$ws.Cells.Item(93, 2).Value = 48
$ws.Cells.Item(93, 3).Value = 1
$ws.Cells.Item(93, 4).Value = 17
$ws.Cells.Item(93, 5).Value = 31
# Row 94
$ws.Cells.Item(94, 1).Value = "Camboya"
$ws.Cells.Item(94, 2).Value = 47
$ws.Cells.Item(94, 3).Value = 10
$ws.Cells.Item(94, 4).Value = 1
$ws.Cells.Item(94, 5).Value = 46
# Row 95
$ws.Cells.Item(95, 1).Value = "Guadalupe"
$ws.Cells.Item(95, 2).Value = 45
$ws.Cells.Item(95, 3).Value = 12
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(95, 5).Value = 45
$ws.Cells.Item(95, 8).Value = 0
# Row 96
$ws.Cells.Item(96, 1).Value = "Azerbaiyan"
$ws.Cells.Item(96, 2).Value = 44
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 7
$ws.Cells.Item(96, 5).Value = 36
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 8).Value = 1
# Row 97
$ws.Cells.Item(97, 1).Value = "Georgia"
$ws.Cells.Item(97, 2).Value = 43
$ws.Cells.Item(97, 3).Value = 3
$ws.Cells.Item(97, 4).Value = 1
$ws.Cells.Item(97, 6).Value = 1
# Row 98
$ws.Cells.Item(98, 1).Value = "Venezuela"
$ws.Cells.Item(98, 2).Value = 42
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 5).Value = 42
# Row 99
$ws.Cells.Item(99, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(99, 3).Value = 11
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 39
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 8).Value = 0
# Row 111
$ws.Cells.Item(111, 6).Value = 1
